# Updated symbol list on Tue Jan 10 20:25:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume 1h (E) columns are stored as literal text in this
# sheet (e.g. "277.75", "0.75%"), not as numbers. Force the text-number format
# before writing so Excel does not silently reinterpret these as numeric/percent
# values, then clear the temporary formatting so styling stays untouched.
$affected = $ws.Range("D2:E47")
$affected.NumberFormat = "@"

$ws.Range("D2").Value = "277.75"
$ws.Range("E2").Value = "0.75%"
$ws.Range("D3").Value = "27.20"
$ws.Range("E3").Value = "1.45%"
$ws.Range("D4").Value = "4.810"
$ws.Range("E4").Value = "-2.22%"
$ws.Range("D5").Value = "0.06368"
$ws.Range("E5").Value = "0.12%"
$ws.Range("D6").Value = "6.959"
$ws.Range("E6").Value = "-0.04%"
$ws.Range("D7").Value = "1.291"
$ws.Range("E7").Value = "-10.55%"
$ws.Range("D8").Value = "0.8768"
$ws.Range("E8").Value = "-1.36%"
$ws.Range("D9").Value = "0.1523"
$ws.Range("E9").Value = "3.38%"
$ws.Range("D10").Value = "0.05067"
$ws.Range("E10").Value = "-2.89%"
$ws.Range("D11").Value = "0.07497"
$ws.Range("E11").Value = "1.13%"
$ws.Range("D12").Value = "0.02966"
$ws.Range("E12").Value = "-5.80%"
$ws.Range("D13").Value = "0.09020"
$ws.Range("E13").Value = "-0.63%"
$ws.Range("D14").Value = "0.001564"
$ws.Range("E14").Value = "-1.71%"
$ws.Range("D15").Value = "0.0006392"
$ws.Range("E15").Value = "0.98%"
$ws.Range("D16").Value = "0.005981"
$ws.Range("E16").Value = "-1.33%"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").Value = "-1.07%"
$ws.Range("D18").Value = "3.305"
$ws.Range("E18").Value = "-1.68%"
$ws.Range("E19").Value = "-0.39%"
$ws.Range("E20").Value = "0.54%"
$ws.Range("E21").Value = "0.07%"
$ws.Range("D22").Value = "3.913"
$ws.Range("E22").Value = "-0.80%"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").Value = "1.16%"
$ws.Range("E24").Value = "-1.08%"
$ws.Range("D25").Value = "0.003871"
$ws.Range("E25").Value = "5.69%"
$ws.Range("D26").Value = "0.0001198"
$ws.Range("E26").Value = "-0.48%"
$ws.Range("E27").Value = "-0.41%"
$ws.Range("D40").Value = "0.04141"
$ws.Range("E40").Value = "2.53%"
$ws.Range("D41").Value = "0.006774"
$ws.Range("E41").Value = "2.22%"
$ws.Range("D42").Value = "0.1178"
$ws.Range("E42").Value = "0.57%"
$ws.Range("D43").Value = "0.002016"
$ws.Range("E43").Value = "-14.82%"
$ws.Range("D44").Value = "0.01120"
$ws.Range("E44").Value = "-8.30%"
$ws.Range("D45").Value = "0.00005167"
$ws.Range("E45").Value = "-1.98%"
$ws.Range("D47").Value = "0.02021"
$ws.Range("E47").Value = "-5.05%"

$affected.ClearFormats()
